# Domino JTAG SPI BOM - add series resistor (R1, 1k) in front of the SPI Flash
# part, renumber the downstream BOM rows and fix up the reference designator
# lists that shifted because of the new part (R1 is now taken by the new
# resistor, everything after it bumps up by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room for the new BOM line: insert a whole row above the old row 8
#    ("10k" resistor bank), shifting everything from row 8 down to row 9.
$ws.Rows.Item(8).Insert()

# 2) Populate the new row 8 with the new 1k series resistor (R1).
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "1k"
$ws.Range("D8").Value = "ANY"
$ws.Range("E8").Value = "R0402_1k_5%_62.5mW"
$ws.Range("F8").Value = "R0402"
$ws.Range("G8").Value = "R1"
$ws.Range("H8").Value = "RES 1K OHM 1/16W 5% 0402 SMD"

# 3) Renumber the "Item" column for the rows pushed down by the insert, and
#    refresh the reference-designator lists that changed because R1 is now
#    used by the new part (10k bank: R1..R4,R6,R7 -> R2..R5,R7,R8; 0R bank:
#    R5,R8,R9,R10..R14 -> R6,R9,R10,R11..R15).
$ws.Range("A9").Value = 8
$ws.Range("G9").Value = "R2, R3, R4, R5, R7, R8"

$ws.Range("A10").Value = 9
$ws.Range("G10").Value = "R6(DNP), R9(DNP), R10(DNP), R11, R12, R13, R14, R15(DNP)"

$ws.Range("A11").Value = 10

$ws.Range("A12").Value = 11

# 4) Extend the print area (and all of the legacy duplicate Print_Area_*
#    defined names the file accumulated) down to the new last row (12), and
#    add one more duplicate entry the same way the original authoring tool
#    did every time the print area was (re)set.
#    The defined names are enumerated in file order; map each position to
#    its new RefersTo value (mirrors the unified diff exactly).
$full12 = "='Domino JTAG SPI Rev. B'!`$A`$1:`$I`$12"
$full1  = "='Domino JTAG SPI Rev. B'!`$A`$1:`$I`$1"

$newRefersTo = @($full12, $full12, $full12, $full12, $full12, $full12, $full12, $full12, $full12, $full12, $full12, $full12, $full1, $full12, $full1, $full1, $full1)

$cnt = $wb.Names.Count()
for ($i = 1; $i -le $cnt; $i++) {
  $n = $wb.Names.Item($i)
  $n.RefersTo = $newRefersTo[$i - 1]
}

$lastPrintAreaName = "_xlnm.Print_Area_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0_0"
$ws.Names.Add($lastPrintAreaName, $full1)

# 5) Match the saved cursor/selection position left behind in the source
#    file (the author's last click before saving).
$ws.Range("C6").Select()
